$d = $word.ActiveDocument

function Replace-InParagraph($index, $oldText, $newText) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $ok = $r.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                           $true, 1, $false, $newText, 2)
    if (-not $ok) {
        throw "Replace failed at paragraph $index for '$oldText'"
    }
}

Replace-InParagraph 63 "Predicate (Attribute): URN. Employee" "Predicate (Attribute): URN. Employs"
Replace-InParagraph 68 "Predicate (Attribute): URN. Employer" "Predicate (Attribute): URN. EmployedBy"
Replace-InParagraph 73 "Predicate (Attribute): URN. Employee" "Predicate (Attribute): URN. PerformedBy"
Replace-InParagraph 78 "Predicate (Attribute): URN. Position" "Predicate (Attribute): URN. Performs"

Replace-InParagraph 84 "Context (Metaclass): SubjectKind. Works, Employs, Performs" `
                       "Context (Metaclass): SubjectKind. Employs, Works, Performs"
Replace-InParagraph 85 "Subject (Context): AggregatedReifiedAggregation SKs. HasEmployer, HasEmployee, HasPosition" `
                       "Subject (Context): AggregatedReifiedAggregation SKs. Employer, Employee, Position"
Replace-InParagraph 86 "Predicate (Role): AggregationSubjectKinds. Employer / Employer, Employee / Position" `
                       "Predicate (Role): AggregationPredicateKinds. Hiring, Work, Performance"
Replace-InParagraph 87 "Object (Occurrence): AggregatedReifiedAggregation OKs. EmployerOf, EmployeeOf, PositionOf" `
                       "Object (Occurrence): AggregatedReifiedAggregation OKs. Employer, Employee, Position"

Replace-InParagraph 91 "Context (Dimension): SubjectKind. Employment" "Context (Dimension): SubjectKind. Employmentship"
Replace-InParagraph 92 "Subject (Measure): AggregatedReifiedActivation SKs. Employed, Employing, Performing" `
                       "Subject (Measure): AggregatedReifiedActivation SKs. Employing, Working, Performing"
Replace-InParagraph 93 "Predicate (Unit): ActivationSubjectKinds. Works, Employs, Performs" `
                       "Predicate (Unit): ActivationPredicateKinds. Employment"
Replace-InParagraph 94 "Object (Value): AggregatedReifiedActivation OKs. EmployeedAt, EmployingWho, PerformingPosition" `
                       "Object (Value): AggregatedReifiedActivation OKs. Employing, Employed, Performing"

Write-Output "done"
